$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 20 (shifts existing rows 20+ down by 2)
$ws.Rows("20:21").Insert()

# Row 20: new weekly price entry
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44608
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100109
$ws.Range("H20").Value = "Uva"
$ws.Range("I20").Value = 100109001
$ws.Range("J20").Value = "Uva"
$ws.Range("K20").Value = "Red Globe"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 10000
$ws.Range("O20").Value = 11000
$ws.Range("P20").Value = 10500
$ws.Range("Q20").Value = "$/bandeja 18 kilos"
$ws.Range("R20").Value = "Región de O'Higgins"
$ws.Range("S20").Value = 583
$ws.Range("T20").Value = 18

# Row 21: new weekly price entry
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44608
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100109
$ws.Range("H21").Value = "Uva"
$ws.Range("I21").Value = 100109001
$ws.Range("J21").Value = "Uva"
$ws.Range("K21").Value = "Superior Seedless"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 9500
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 9750
$ws.Range("Q21").Value = "$/bandeja 18 kilos"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 542
$ws.Range("T21").Value = 18
